# Auto-generated: scheduled market-data refresh for Bahamut_Profits workbook
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for the
# specific leve rows whose underlying market data changed in this run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 394.88  # H28
$ws.Cells.Item(28, 9).Value = 461.1579  # I28
$ws.Cells.Item(28, 10).Value = 185  # J28
$ws.Cells.Item(28, 11).Value = 461.1579  # K28
$ws.Cells.Item(28, 12).Value = 185  # L28
$ws.Cells.Item(28, 13).Value = 23.84210000000002  # M28
$ws.Cells.Item(28, 14).Value = -1155  # N28
$ws.Cells.Item(33, 8).Value = 53057.684  # H33
$ws.Cells.Item(33, 9).Value = 83745.5  # I33
$ws.Cells.Item(33, 10).Value = 450  # J33
$ws.Cells.Item(33, 11).Value = 83745.5  # K33
$ws.Cells.Item(33, 12).Value = 450  # L33
$ws.Cells.Item(33, 13).Value = -83516.5  # M33
$ws.Cells.Item(33, 14).Value = -908  # N33
$ws.Cells.Item(101, 8).Value = 144869.86  # H101
$ws.Cells.Item(101, 9).Value = 2420.8  # I101
$ws.Cells.Item(101, 10).Value = 500992.5  # J101
$ws.Cells.Item(101, 11).Value = 7262.400000000001  # K101
$ws.Cells.Item(101, 12).Value = 1502977.5  # L101
$ws.Cells.Item(101, 13).Value = -5640.400000000001  # M101
$ws.Cells.Item(101, 14).Value = -1506221.5  # N101
$ws.Cells.Item(113, 9).Value = 4926.25  # I113
$ws.Cells.Item(113, 10).Value = 3157.6  # J113
$ws.Cells.Item(113, 11).Value = 4926.25  # K113
$ws.Cells.Item(113, 12).Value = 3157.6  # L113
$ws.Cells.Item(113, 13).Value = -1672.25  # M113
$ws.Cells.Item(113, 14).Value = -9665.6  # N113
$ws.Cells.Item(139, 8).Value = 66383.336  # H139
$ws.Cells.Item(139, 10).Value = 66383.336  # J139
$ws.Cells.Item(139, 12).Value = 66383.336  # L139
$ws.Cells.Item(139, 14).Value = -76663.336  # N139

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1521.8235  # H61
$ws.Cells.Item(61, 9).Value = 1491.4667  # I61
$ws.Cells.Item(61, 10).Value = 1749.5  # J61
$ws.Cells.Item(61, 11).Value = 1491.4667  # K61
$ws.Cells.Item(61, 12).Value = 1749.5  # L61
$ws.Cells.Item(61, 13).Value = -1279.4667  # M61
$ws.Cells.Item(61, 14).Value = -2173.5  # N61
$ws.Cells.Item(70, 8).Value = 0  # H70
$ws.Cells.Item(70, 10).Value = 0  # J70
$ws.Cells.Item(70, 12).Value = 0  # L70
$ws.Cells.Item(70, 14).ClearContents()  # N70
$ws.Cells.Item(73, 8).Value = 0  # H73
$ws.Cells.Item(73, 10).Value = 0  # J73
$ws.Cells.Item(73, 12).Value = 0  # L73
$ws.Cells.Item(73, 14).ClearContents()  # N73
$ws.Cells.Item(136, 8).Value = 1521.8235  # H136
$ws.Cells.Item(136, 9).Value = 1491.4667  # I136
$ws.Cells.Item(136, 10).Value = 1749.5  # J136
$ws.Cells.Item(136, 11).Value = 4474.4001  # K136
$ws.Cells.Item(136, 12).Value = 5248.5  # L136
$ws.Cells.Item(136, 13).Value = -1924.4001  # M136
$ws.Cells.Item(136, 14).Value = -10348.5  # N136
$ws.Cells.Item(139, 8).Value = 34322.6  # H139
$ws.Cells.Item(139, 10).Value = 34322.6  # J139
$ws.Cells.Item(139, 12).Value = 34322.6  # L139
$ws.Cells.Item(139, 14).Value = -44602.6  # N139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 1051.75  # H80
$ws.Cells.Item(80, 9).Value = 750  # I80
$ws.Cells.Item(80, 10).Value = 1152.3334  # J80
$ws.Cells.Item(80, 11).Value = 750  # K80
$ws.Cells.Item(80, 12).Value = 1152.3334  # L80
$ws.Cells.Item(80, 13).Value = 248  # M80
$ws.Cells.Item(80, 14).Value = -3148.3334  # N80
$ws.Cells.Item(83, 8).Value = 1051.75  # H83
$ws.Cells.Item(83, 9).Value = 750  # I83
$ws.Cells.Item(83, 10).Value = 1152.3334  # J83
$ws.Cells.Item(83, 11).Value = 3750  # K83
$ws.Cells.Item(83, 12).Value = 5761.666999999999  # L83
$ws.Cells.Item(83, 13).Value = 1242  # M83
$ws.Cells.Item(83, 14).Value = -15745.667  # N83
$ws.Cells.Item(107, 8).Value = 9199.294  # H107
$ws.Cells.Item(107, 9).Value = 1479.2  # I107
$ws.Cells.Item(107, 10).Value = 67100  # J107
$ws.Cells.Item(107, 11).Value = 1479.2  # K107
$ws.Cells.Item(107, 12).Value = 67100  # L107
$ws.Cells.Item(107, 13).Value = 440.8  # M107
$ws.Cells.Item(107, 14).Value = -70940  # N107
$ws.Cells.Item(134, 8).Value = 41766.88  # H134
$ws.Cells.Item(134, 9).Value = 1689.6364  # I134
$ws.Cells.Item(134, 10).Value = 335666.66  # J134
$ws.Cells.Item(134, 11).Value = 5068.9092  # K134
$ws.Cells.Item(134, 12).Value = 1006999.98  # L134
$ws.Cells.Item(134, 13).Value = -2533.9092  # M134
$ws.Cells.Item(134, 14).Value = -1012069.98  # N134
$ws.Cells.Item(137, 8).Value = 40587.832  # H137
$ws.Cells.Item(137, 10).Value = 40587.832  # J137
$ws.Cells.Item(137, 12).Value = 40587.832  # L137
$ws.Cells.Item(137, 14).Value = -50787.832  # N137

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1649.4828  # H31
$ws.Cells.Item(31, 9).Value = 1416.75  # I31
$ws.Cells.Item(31, 10).Value = 2166.6667  # J31
$ws.Cells.Item(31, 11).Value = 1416.75  # K31
$ws.Cells.Item(31, 12).Value = 2166.6667  # L31
$ws.Cells.Item(31, 13).Value = -1121.75  # M31
$ws.Cells.Item(31, 14).Value = -2756.6667  # N31
$ws.Cells.Item(34, 8).Value = 1649.4828  # H34
$ws.Cells.Item(34, 9).Value = 1416.75  # I34
$ws.Cells.Item(34, 10).Value = 2166.6667  # J34
$ws.Cells.Item(34, 11).Value = 1416.75  # K34
$ws.Cells.Item(34, 12).Value = 2166.6667  # L34
$ws.Cells.Item(34, 13).Value = -1214.75  # M34
$ws.Cells.Item(34, 14).Value = -2570.6667  # N34
$ws.Cells.Item(92, 8).Value = 0  # H92
$ws.Cells.Item(92, 10).Value = 0  # J92
$ws.Cells.Item(92, 12).Value = 0  # L92
$ws.Cells.Item(92, 14).ClearContents()  # N92
$ws.Cells.Item(99, 8).Value = 2782.5652  # H99
$ws.Cells.Item(99, 9).Value = 2486.6667  # I99
$ws.Cells.Item(99, 11).Value = 2486.6667  # K99
$ws.Cells.Item(99, 13).Value = -988.6667000000002  # M99
$ws.Cells.Item(126, 8).Value = 2782.5652  # H126
$ws.Cells.Item(126, 9).Value = 2486.6667  # I126
$ws.Cells.Item(126, 11).Value = 7460.000100000001  # K126
$ws.Cells.Item(126, 13).Value = -4990.000100000001  # M126

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 300  # H92
$ws.Cells.Item(92, 9).Value = 350  # I92
$ws.Cells.Item(92, 10).Value = 200  # J92
$ws.Cells.Item(92, 11).Value = 1050  # K92
$ws.Cells.Item(92, 12).Value = 600  # L92
$ws.Cells.Item(92, 13).Value = 198  # M92
$ws.Cells.Item(92, 14).Value = -3096  # N92
$ws.Cells.Item(109, 8).Value = 2438.7778  # H109
$ws.Cells.Item(109, 9).Value = 779.6  # I109
$ws.Cells.Item(109, 11).Value = 2338.8  # K109
$ws.Cells.Item(109, 13).Value = -1298.8  # M109
$ws.Cells.Item(112, 8).Value = 3710.6667  # H112
$ws.Cells.Item(112, 9).Value = 940  # I112
$ws.Cells.Item(112, 10).Value = 4136.923  # J112
$ws.Cells.Item(112, 11).Value = 2820  # K112
$ws.Cells.Item(112, 12).Value = 12410.769  # L112
$ws.Cells.Item(112, 13).Value = -1712  # M112
$ws.Cells.Item(112, 14).Value = -14626.769  # N112

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4811.9  # H70
$ws.Cells.Item(70, 9).Value = 4077.375  # I70
$ws.Cells.Item(70, 11).Value = 4077.375  # K70
$ws.Cells.Item(70, 13).Value = -3807.375  # M70
$ws.Cells.Item(73, 8).Value = 4811.9  # H73
$ws.Cells.Item(73, 9).Value = 4077.375  # I73
$ws.Cells.Item(73, 11).Value = 4077.375  # K73
$ws.Cells.Item(73, 13).Value = -3141.375  # M73
$ws.Cells.Item(94, 8).Value = 28000  # H94
$ws.Cells.Item(94, 10).Value = 28000  # J94
$ws.Cells.Item(94, 12).Value = 28000  # L94
$ws.Cells.Item(94, 14).Value = -29352  # N94
$ws.Cells.Item(136, 8).Value = 19430.334  # H136
$ws.Cells.Item(136, 10).Value = 19430.334  # J136
$ws.Cells.Item(136, 12).Value = 58291.00199999999  # L136
$ws.Cells.Item(136, 14).Value = -63391.00199999999  # N136
$ws.Cells.Item(139, 8).Value = 44577  # H139
$ws.Cells.Item(139, 10).Value = 44577  # J139
$ws.Cells.Item(139, 12).Value = 44577  # L139
$ws.Cells.Item(139, 14).Value = -54857  # N139

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 412.27274  # H16
$ws.Cells.Item(16, 9).Value = 417.61905  # I16
$ws.Cells.Item(16, 11).Value = 417.61905  # K16
$ws.Cells.Item(16, 13).Value = -247.61905  # M16
$ws.Cells.Item(93, 8).Value = 1249.95  # H93
$ws.Cells.Item(93, 9).Value = 573.6667  # I93
$ws.Cells.Item(93, 10).Value = 3278.8  # J93
$ws.Cells.Item(93, 11).Value = 573.6667  # K93
$ws.Cells.Item(93, 12).Value = 3278.8  # L93
$ws.Cells.Item(93, 13).Value = 674.3333  # M93
$ws.Cells.Item(93, 14).Value = -5774.8  # N93
$ws.Cells.Item(134, 8).Value = 38582.25  # H134
$ws.Cells.Item(134, 10).Value = 38582.25  # J134
$ws.Cells.Item(134, 12).Value = 38582.25  # L134
$ws.Cells.Item(134, 14).Value = -48722.25  # N134
$ws.Cells.Item(138, 8).Value = 35266.668  # H138
$ws.Cells.Item(138, 10).Value = 35266.668  # J138
$ws.Cells.Item(138, 12).Value = 35266.668  # L138
$ws.Cells.Item(138, 14).Value = -45546.668  # N138

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(138, 8).Value = 46932  # H138
$ws.Cells.Item(138, 10).Value = 46932  # J138
$ws.Cells.Item(138, 12).Value = 46932  # L138
$ws.Cells.Item(138, 14).Value = -57212  # N138
